$d = $word.ActiveDocument

$replacements = @(
    @('2024-05-11 Saturday', '2024-05-12 Sunday'),
    @('72×50=3600', '27×29=783'),
    @('43×57=2451', '24×65=1560'),
    @('42×46=1932', '42×66=2772'),
    @('68×12=816', '73×86=6278'),
    @('93×36=3348', '70×20=1400'),
    @('46×89=4094', '73×87=6351'),
    @('97×59=5723', '11×73=803'),
    @('63×19=1197', '40×28=1120'),
    @('35×68=2380', '77×60=4620'),
    @('24×53=1272', '43×27=1161'),
    @('27×73=1971', '75×28=2100'),
    @('20×97=1940', '71×18=1278'),
    @('16×60=960', '92×56=5152'),
    @('44×60=2640', '72×22=1584'),
    @('47×93=4371', '12×91=1092'),
    @('44×48=2112', '75×78=5850'),
    @('69×92=6348', '82×21=1722'),
    @('16×22=352', '81×93=7533'),
    @('39×21=819', '14×27=378'),
    @('57×50=2850', '96×23=2208'),
    @('78×12=936', '86×75=6450'),
    @('91×77=7007', '22×14=308'),
    @('17×30=510', '75×44=3300'),
    @('46×41=1886', '72×22=1584'),
    @('44×95=4180', '47×96=4512'),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
